# Updated cryptos list -- refresh Price (D) and Volume(1h) (E) columns, and
# restore the two rank swaps (PancakeSwap/NEARProtocol and Kaspa/dogwifhat)
# produced by this run's scrape.
#
# Cells whose new text looks like a plain number (e.g. "0.994", "1.00") are
# written with the column's NumberFormat temporarily forced to Text ("@") so
# Excel keeps the literal digits (incl. trailing zeros) instead of coercing
# them to a float; the cell style is then reset to "Normal" so no stray
# formatting delta is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "68.066.54"; ForceText = $false },
    @{ Cell = "E2"; Value = "  +0.51%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "3.790.71"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -0.35%  "; ForceText = $false },
    @{ Cell = "D4"; Value = "0.994"; ForceText = $true },
    @{ Cell = "E4"; Value = "  -0.56%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "600.36"; ForceText = $true },
    @{ Cell = "E5"; Value = "  +0.64%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "164.71"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -1.46%  "; ForceText = $false },
    @{ Cell = "E7"; Value = "  -0.06%  "; ForceText = $false },
    @{ Cell = "E8"; Value = "  -0.68%  "; ForceText = $false },
    @{ Cell = "D9"; Value = "0.158"; ForceText = $true },
    @{ Cell = "E9"; Value = "  -1.43%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "0.449"; ForceText = $true },
    @{ Cell = "E10"; Value = "  +0.01%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "6.51"; ForceText = $true },
    @{ Cell = "E11"; Value = "  +3.34%  "; ForceText = $false },
    @{ Cell = "E12"; Value = "  -2.54%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "35.60"; ForceText = $true },
    @{ Cell = "E13"; Value = "  -0.93%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "4.427.90"; ForceText = $false },
    @{ Cell = "E14"; Value = "  -0.73%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "3.838.17"; ForceText = $false },
    @{ Cell = "E15"; Value = "  +0.62%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "68.050.32"; ForceText = $false },
    @{ Cell = "E16"; Value = "  +0.54%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "18.28"; ForceText = $true },
    @{ Cell = "E17"; Value = "  -1.55%  "; ForceText = $false },
    @{ Cell = "E18"; Value = "  +2.19%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "7.05"; ForceText = $true },
    @{ Cell = "E19"; Value = "  -0.56%  "; ForceText = $false },
    @{ Cell = "D20"; Value = "460.65"; ForceText = $true },
    @{ Cell = "E20"; Value = "  -0.08%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "9.67"; ForceText = $true },
    @{ Cell = "E21"; Value = "  -2.68%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "0.700"; ForceText = $true },
    @{ Cell = "E22"; Value = "  -0.21%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "0.0000146"; ForceText = $true },
    @{ Cell = "E23"; Value = "  -4.76%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "83.01"; ForceText = $true },
    @{ Cell = "E24"; Value = "  -0.67%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "11.96"; ForceText = $true },
    @{ Cell = "E25"; Value = "  -1.19%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "2.10"; ForceText = $true },
    @{ Cell = "E26"; Value = "  -0.42%  "; ForceText = $false },
    @{ Cell = "E27"; Value = "  -0.09%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "9.96"; ForceText = $true },
    @{ Cell = "E28"; Value = "  -0.57%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "3.942.18"; ForceText = $false },
    @{ Cell = "E29"; Value = "  -0.33%  "; ForceText = $false },
    @{ Cell = "E30"; Value = "  -0.24%  "; ForceText = $false },
    @{ Cell = "B31"; Value = "PancakeSwap"; ForceText = $false },
    @{ Cell = "C31"; Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; ForceText = $false },
    @{ Cell = "D31"; Value = "2.63"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -4.95%  "; ForceText = $false },
    @{ Cell = "B32"; Value = "NEARProtocol"; ForceText = $false },
    @{ Cell = "C32"; Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; ForceText = $false },
    @{ Cell = "D32"; Value = "7.32"; ForceText = $true },
    @{ Cell = "E32"; Value = "  +0.32%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "29.23"; ForceText = $true },
    @{ Cell = "E33"; Value = "  -1.22%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E34"; Value = "  +0.09%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "9.00"; ForceText = $true },
    @{ Cell = "E35"; Value = "  -1.06%  "; ForceText = $false },
    @{ Cell = "D36"; Value = "0.0996"; ForceText = $true },
    @{ Cell = "E36"; Value = "  -0.46%  "; ForceText = $false },
    @{ Cell = "B37"; Value = "Kaspa"; ForceText = $false },
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; ForceText = $false },
    @{ Cell = "D37"; Value = "0.139"; ForceText = $true },
    @{ Cell = "E37"; Value = "  +1.07%  "; ForceText = $false },
    @{ Cell = "B38"; Value = "dogwifhat"; ForceText = $false },
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; ForceText = $false },
    @{ Cell = "D38"; Value = "3.30"; ForceText = $true },
    @{ Cell = "E38"; Value = "  -2.84%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "5.81"; ForceText = $true },
    @{ Cell = "E39"; Value = "  +0.40%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "0.987"; ForceText = $true },
    @{ Cell = "E40"; Value = "  -1.36%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "0.999"; ForceText = $true },
    @{ Cell = "E41"; Value = "  +0.04%  "; ForceText = $false },
    @{ Cell = "E42"; Value = "  +0.02%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "47.52"; ForceText = $true },
    @{ Cell = "E43"; Value = "  -1.41%  "; ForceText = $false },
    @{ Cell = "D44"; Value = "0.299"; ForceText = $true },
    @{ Cell = "E44"; Value = "  -0.45%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "43.19"; ForceText = $true },
    @{ Cell = "E45"; Value = "  -0.08%  "; ForceText = $false },
    @{ Cell = "D46"; Value = "151.93"; ForceText = $true },
    @{ Cell = "E46"; Value = "  +2.55%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "8.35"; ForceText = $true },
    @{ Cell = "E47"; Value = "  +0.30%  "; ForceText = $false },
    @{ Cell = "E48"; Value = "  +1.59%  "; ForceText = $false },
    @{ Cell = "E49"; Value = "  +1.24%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "389.10"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -1.34%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "26.52"; ForceText = $true },
    @{ Cell = "E51"; Value = "  -1.25%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
